# The sheet contains a weekly table of "Haba" (fava bean) price records for
# "Femacal de La Calera" in "Coquimbo". A new weekly record needs to be
# inserted as row 146 (pushing the existing rows 146:150 down to 147:151),
# growing the used range from A1:R150 to A1:R151.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 146, shifting rows 146:150 down to 147:151.
$ws.Rows("146:146").Insert()

# Populate the new row 146 with the new weekly price record.
$ws.Cells.Item(146, 1).Value  = 3
$ws.Cells.Item(146, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(146, 3).Value  = "Coquimbo"
$ws.Cells.Item(146, 4).Value  = 44746
$ws.Cells.Item(146, 5).Value  = 5
$ws.Cells.Item(146, 6).Value  = 100112026
$ws.Cells.Item(146, 7).Value  = "Haba"
$ws.Cells.Item(146, 8).Value  = "Sin especificar"
$ws.Cells.Item(146, 9).Value  = "Primera"
$ws.Cells.Item(146, 10).Value = 101
$ws.Cells.Item(146, 11).Value = 16000
$ws.Cells.Item(146, 12).Value = 17000
$ws.Cells.Item(146, 13).Value = 16446
$ws.Cells.Item(146, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(146, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(146, 16).Value = 658
$ws.Cells.Item(146, 17).Value = 25
$ws.Cells.Item(146, 18).Value = "Hortaliza"
